$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034873433281631
$ws.Range("D2").Value = 1.053126138092693
$ws.Range("E2").Value = 1.045281254983431
$ws.Range("F2").Value = 1.059286841184261
$ws.Range("I2").Value = 1.045630254507131
$ws.Range("J2").Value = 1.039990295279646
$ws.Range("K2").Value = 1.05587296836894
$ws.Range("L2").Value = 1.048049935467291
$ws.Range("M2").Value = 1.062016757899615
$ws.Range("N2").Value = 1.017296294447743
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035650538168227
$ws.Range("D3").Value = 1.053698481775425
$ws.Range("E3").Value = 1.045950196853032
$ws.Range("F3").Value = 1.059986208680758
$ws.Range("I3").Value = 1.045814529678178
$ws.Range("J3").Value = 1.040411584607797
$ws.Range("K3").Value = 1.056258911672309
$ws.Range("L3").Value = 1.048530657254337
$ws.Range("M3").Value = 1.062530615879904
$ws.Range("N3").Value = 1.01743658878861
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036154113877209
$ws.Range("D4").Value = 1.05406924644732
$ws.Range("E4").Value = 1.046384032282273
$ws.Range("F4").Value = 1.060439627131061
$ws.Range("I4").Value = 1.045932758069468
$ws.Range("J4").Value = 1.040684233004692
$ws.Range("K4").Value = 1.056508358604513
$ws.Range("L4").Value = 1.048842017381493
$ws.Range("M4").Value = 1.062863318091853
$ws.Range("N4").Value = 1.017527360793525
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036365991581202
$ws.Range("D5").Value = 1.054225214399474
$ws.Range("E5").Value = 1.046566650516869
$ws.Range("F5").Value = 1.060630452914912
$ws.Range("I5").Value = 1.045982218544155
$ws.Range("J5").Value = 1.040798864042853
$ws.Range("K5").Value = 1.056613156530581
$ws.Range("L5").Value = 1.048972983664475
$ws.Range("M5").Value = 1.063003232627193
$ws.Range("N5").Value = 1.017565518985633
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03640157695511
$ws.Range("D6").Value = 1.054251407819691
$ws.Range("E6").Value = 1.046597326573232
$ws.Range("F6").Value = 1.06066250557526
$ws.Range("I6").Value = 1.045990508913341
$ws.Range("J6").Value = 1.040818111624773
$ws.Range("K6").Value = 1.056630748439745
$ws.Range("L6").Value = 1.048994977567246
$ws.Range("M6").Value = 1.063026727570055
$ws.Range("N6").Value = 1.017571925761475
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036156944314315
$ws.Range("D7").Value = 1.054071330114371
$ws.Range("E7").Value = 1.04638647151991
$ws.Range("F7").Value = 1.06044217613741
$ws.Range("I7").Value = 1.045933419917751
$ws.Range("J7").Value = 1.040685764674027
$ws.Range("K7").Value = 1.056509759194793
$ws.Range("L7").Value = 1.048843767084393
$ws.Range("M7").Value = 1.06286518745614
$ws.Range("N7").Value = 1.017527870675037
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035135905775697
$ws.Range("D8").Value = 1.053319476102978
$ws.Range("E8").Value = 1.045507121765188
$ws.Range("F8").Value = 1.059523011968039
$ws.Range("I8").Value = 1.045692739711148
$ws.Range("J8").Value = 1.040132661627579
$ws.Range("K8").Value = 1.056003457641202
$ws.Range("L8").Value = 1.04821233438082
$ws.Range("M8").Value = 1.062190375774698
$ws.Range("N8").Value = 1.01734370890208
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033342434749367
$ws.Range("D9").Value = 1.05199792185546
$ws.Range("E9").Value = 1.043965232683184
$ws.Range("F9").Value = 1.0579101723655
$ws.Range("I9").Value = 1.045260937468649
$ws.Range("J9").Value = 1.039158436407627
$ws.Range("K9").Value = 1.055109180564346
$ws.Range("L9").Value = 1.04710204391964
$ws.Range("M9").Value = 1.061002894719243
$ws.Range("N9").Value = 1.017019154037091
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032150755245512
$ws.Range("D10").Value = 1.051119234131028
$ws.Range("E10").Value = 1.042942565156282
$ws.Range("F10").Value = 1.056839681305418
$ws.Range("I10").Value = 1.044967954000185
$ws.Range("J10").Value = 1.038509312372551
$ws.Range("K10").Value = 1.054511669165591
$ws.Range("L10").Value = 1.046363540779519
$ws.Range("M10").Value = 1.060212438008811
$ws.Range("N10").Value = 1.016802787530996
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031635710483261
$ws.Range("D11").Value = 1.050739335961364
$ws.Range("E11").Value = 1.042501012286573
$ws.Range("F11").Value = 1.056377298259876
$ws.Range("I11").Value = 1.044839887701629
$ws.Range("J11").Value = 1.03822833707079
$ws.Range("K11").Value = 1.054252644817782
$ws.Range("L11").Value = 1.046044181644604
$ws.Range("M11").Value = 1.059870467967196
$ws.Range("N11").Value = 1.016709105426819
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031444546190811
$ws.Range("D12").Value = 1.050598314060281
$ws.Range("E12").Value = 1.042337192673179
$ws.Range("F12").Value = 1.056205723185945
$ws.Range("I12").Value = 1.04479213847387
$ws.Range("J12").Value = 1.038123986713566
$ws.Range("K12").Value = 1.054156388443018
$ws.Range("L12").Value = 1.045925621801379
$ws.Range("M12").Value = 1.059743492404449
$ws.Range("N12").Value = 1.016674309127093
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031485544929779
$ws.Range("D13").Value = 1.050628559666192
$ws.Range("E13").Value = 1.04237232377819
$ws.Range("F13").Value = 1.056242518684787
$ws.Range("I13").Value = 1.044802388961754
$ws.Range("J13").Value = 1.038146369455661
$ws.Range("K13").Value = 1.054177037691858
$ws.Range("L13").Value = 1.045951050351858
$ws.Range("M13").Value = 1.059770726922322
$ws.Range("N13").Value = 1.016681772979899
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031619905787967
$ws.Range("D14").Value = 1.050727677215119
$ws.Range("E14").Value = 1.042487466961149
$ws.Range("F14").Value = 1.05636311223611
$ws.Range("I14").Value = 1.044835944395868
$ws.Range("J14").Value = 1.038219711095447
$ws.Range("K14").Value = 1.054244689108007
$ws.Range("L14").Value = 1.046034380128719
$ws.Range("M14").Value = 1.059859971151366
$ws.Range("N14").Value = 1.016706229122277
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031702709388082
$ws.Range("D15").Value = 1.050788758691412
$ws.Range("E15").Value = 1.042558436085926
$ws.Range("F15").Value = 1.056437437105179
$ws.Range("I15").Value = 1.044856595222811
$ws.Range("J15").Value = 1.038264901512791
$ws.Range("K15").Value = 1.054286365700837
$ws.Range("L15").Value = 1.046085730931055
$ws.Range("M15").Value = 1.059914963793449
$ws.Range("N15").Value = 1.016721297559875
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032184957472636
$ws.Range("D16").Value = 1.051144459112082
$ws.Range("E16").Value = 1.042971896483611
$ws.Range("F16").Value = 1.05687039250082
$ws.Range("I16").Value = 1.044976428085071
$ws.Range("J16").Value = 1.038527961998094
$ws.Range("K16").Value = 1.054528853614482
$ws.Range("L16").Value = 1.046384744530136
$ws.Range("M16").Value = 1.060235139972476
$ws.Range("N16").Value = 1.016809005073299
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032487717463879
$ws.Range("D17").Value = 1.051367737229048
$ws.Range("E17").Value = 1.043231590506462
$ws.Range("F17").Value = 1.057142282601421
$ws.Range("I17").Value = 1.045051274803829
$ws.Range("J17").Value = 1.038693000571079
$ws.Range("K17").Value = 1.054680881208912
$ws.Range("L17").Value = 1.046572420807986
$ws.Range("M17").Value = 1.06043606042304
$ws.Range("N17").Value = 1.016864023644701
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032664404768817
$ws.Range("D18").Value = 1.051498027338641
$ws.Range("E18").Value = 1.043383187925097
$ws.Range("F18").Value = 1.057300981977565
$ws.Range("I18").Value = 1.045094815537899
$ws.Range("J18").Value = 1.03878927430299
$ws.Range("K18").Value = 1.054769527477958
$ws.Range("L18").Value = 1.046681929389319
$ws.Range("M18").Value = 1.060553283028697
$ws.Range("N18").Value = 1.016896115593959
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03272466620711
$ws.Range("D19").Value = 1.051542462291835
$ws.Range("E19").Value = 1.043434899404291
$ws.Range("F19").Value = 1.057355113032855
$ws.Range("I19").Value = 1.045109642082731
$ws.Range("J19").Value = 1.038822102730451
$ws.Range("K19").Value = 1.054799748622892
$ws.Range("L19").Value = 1.04671927575214
$ws.Range("M19").Value = 1.06059325778926
$ws.Range("N19").Value = 1.016907058196939
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032455224597564
$ws.Range("D20").Value = 1.051343775817042
$ws.Range("E20").Value = 1.04320371513118
$ws.Range("F20").Value = 1.057113099922626
$ws.Range("I20").Value = 1.045043256460109
$ws.Range("J20").Value = 1.038675292498823
$ws.Range("K20").Value = 1.054664573060075
$ws.Range("L20").Value = 1.046552280751456
$ws.Range("M20").Value = 1.060414500536439
$ws.Range("N20").Value = 1.016858120612451
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03158033581051
$ws.Range("D21").Value = 1.050698487088018
$ws.Range("E21").Value = 1.042453554824715
$ws.Range("F21").Value = 1.056327595603292
$ws.Range("I21").Value = 1.044826068109796
$ws.Range("J21").Value = 1.038198113347086
$ws.Range("K21").Value = 1.054224768632045
$ws.Range("L21").Value = 1.046009839796759
$ws.Range("M21").Value = 1.059833689628896
$ws.Range("N21").Value = 1.016699027355102
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031031105731914
$ws.Range("D22").Value = 1.050293285559191
$ws.Range("E22").Value = 1.041983015703657
$ws.Range("F22").Value = 1.055834729046555
$ws.Range("I22").Value = 1.044688473948572
$ws.Range("J22").Value = 1.037898187339941
$ws.Range("K22").Value = 1.053947996970582
$ws.Range("L22").Value = 1.045669158799887
$ws.Range("M22").Value = 1.059468785536916
$ws.Range("N22").Value = 1.016599007485599
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03132218201753
$ws.Range("D23").Value = 1.050508040826225
$ws.Range("E23").Value = 1.042232350855211
$ws.Range("F23").Value = 1.056095910316201
$ws.Range("I23").Value = 1.044761513401627
$ws.Range("J23").Value = 1.038057174274006
$ws.Range("K23").Value = 1.054094741959972
$ws.Range("L23").Value = 1.0458497243308
$ws.Range("M23").Value = 1.059662201445944
$ws.Range("N23").Value = 1.01665202895005
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03246990643959
$ws.Range("D24").Value = 1.051354602775719
$ws.Range("E24").Value = 1.04321631043583
$ws.Range("F24").Value = 1.057126285977968
$ws.Range("I24").Value = 1.045046879963318
$ws.Range("J24").Value = 1.038683293985199
$ws.Range("K24").Value = 1.054671942099925
$ws.Range("L24").Value = 1.046561381052141
$ws.Range("M24").Value = 1.060424242431472
$ws.Range("N24").Value = 1.016860787937007
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033805399119424
$ws.Range("D25").Value = 1.052339170432223
$ws.Range("E25").Value = 1.044362930144119
$ws.Range("F25").Value = 1.058326305153927
$ws.Range("I25").Value = 1.045373474482049
$ws.Range("J25").Value = 1.039410239605817
$ws.Range("K25").Value = 1.05534061313832
$ws.Range("L25").Value = 1.047388789523493
$ws.Range("M25").Value = 1.061309683740412
$ws.Range("N25").Value = 1.017103060780794
